$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26; this shifts the existing rows 26-28
# down to 27-29 (their content, including the D-column date style,
# travels with them) and leaves a fresh blank row 26 behind.
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new weekly record.
$ws.Cells.Item(26, 1).Value2 = 5
$ws.Cells.Item(26, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(26, 3).Value = "Maule"
$ws.Cells.Item(26, 4).Value2 = 44522
$ws.Cells.Item(26, 5).Value2 = 7
$ws.Cells.Item(26, 6).Value = "Fruta"
$ws.Cells.Item(26, 7).Value2 = 100101
$ws.Cells.Item(26, 8).Value = "Berries"
$ws.Cells.Item(26, 9).Value2 = 100101001
$ws.Cells.Item(26, 10).Value = "Arándano (blue)"
$ws.Cells.Item(26, 11).Value = "Sin especificar"
$ws.Cells.Item(26, 12).Value = "Primera"
$ws.Cells.Item(26, 13).Value2 = 30
$ws.Cells.Item(26, 14).Value2 = 5000
$ws.Cells.Item(26, 15).Value2 = 5000
$ws.Cells.Item(26, 16).Value2 = 5000
$ws.Cells.Item(26, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(26, 18).Value = "Provincia de Linares"
$ws.Cells.Item(26, 19).Value2 = 2500
$ws.Cells.Item(26, 20).Value2 = 2
